$d = $word.ActiveDocument

# The commit adds explicit single-line-spacing (w:line="240" w:lineRule="auto")
# paragraph formatting to all of the "body" style paragraphs (Heading 5,
# First Paragraph, Compact and Body Text), while leaving Title/Author/Date/
# Abstract/Heading1-4/SourceCode paragraphs untouched.
foreach ($p in $d.Paragraphs) {
    $styleName = $p.Style.NameLocal
    if ($styleName -eq "Heading 5" -or `
        $styleName -eq "First Paragraph" -or `
        $styleName -eq "Compact" -or `
        $styleName -eq "Body Text") {
        $p.Range.ParagraphFormat.LineSpacingRule = 0   # wdLineSpaceSingle -> line=240 auto
    }
}

# The "Body Text" style itself also picks up an explicit double line spacing
# (w:line="480" w:lineRule="auto") in its style definition, in addition to
# the existing before/after spacing.
$bodyTextStyle = $d.Styles.Item("Body Text")
$bodyTextStyle.ParagraphFormat.LineSpacingRule = 2   # wdLineSpaceDouble -> line=480 auto
